$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.336.92'
$ws.Range("E2").Value = '  +0.57%  '
$ws.Range("D3").Value = '2.564.43'
$ws.Range("E3").Value = '  +0.20%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '582.73'
$ws.Range("E5").Value = '  +0.31%  '
$ws.Range("D6").Value = '171.71'
$ws.Range("E6").Value = '  +0.75%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  +1.44%  '
$ws.Range("D9").Value = '2.562.63'
$ws.Range("E9").Value = '  +0.16%  '
$ws.Range("D10").Value = '0.166'
$ws.Range("E10").Value = '  -0.49%  '
$ws.Range("E11").Value = '  -0.48%  '
$ws.Range("D12").Value = '0.360'
$ws.Range("E12").Value = '  +2.71%  '
$ws.Range("E13").Value = '  +1.79%  '
$ws.Range("D14").Value = '3.040.96'
$ws.Range("E14").Value = '  -0.46%  '
$ws.Range("D15").Value = '71.196.77'
$ws.Range("E15").Value = '  +0.53%  '
$ws.Range("D16").Value = '0.0000180'
$ws.Range("E16").Value = '  -2.88%  '
$ws.Range("D17").Value = '25.49'
$ws.Range("E17").Value = '  +0.95%  '
$ws.Range("D18").Value = '2.580.32'
$ws.Range("E18").Value = '  +0.24%  '
$ws.Range("D19").Value = '11.66'
$ws.Range("E19").Value = '  -1.90%  '
$ws.Range("E20").Value = '  +3.58%  '
$ws.Range("D21").Value = '357.81'
$ws.Range("E21").Value = '  -1.88%  '
$ws.Range("D22").Value = '3.97'
$ws.Range("E22").Value = '  -0.93%  '
$ws.Range("D23").Value = '2.05'
$ws.Range("E23").Value = '  +3.12%  '
$ws.Range("E24").Value = '  +0.08%  '
$ws.Range("D25").Value = '70.60'
$ws.Range("E25").Value = '  +0.35%  '
$ws.Range("D26").Value = '4.11'
$ws.Range("E26").Value = '  -1.03%  '
$ws.Range("D27").Value = '9.15'
$ws.Range("E27").Value = '  -1.34%  '
$ws.Range("E29").Value = '  +0.36%  '
$ws.Range("D30").Value = '0.0₃0925'
$ws.Range("E30").Value = '  -0.46%  '
$ws.Range("D31").Value = '7.99'
$ws.Range("E31").Value = '  +2.33%  '
$ws.Range("D32").Value = '475.89'
$ws.Range("E32").Value = '  -2.09%  '
$ws.Range("E33").Value = '  -1.20%  '
$ws.Range("D34").Value = '1.77'
$ws.Range("E34").Value = '  -0.15%  '
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("E36").Value = '  +3.68%  '
$ws.Range("D37").Value = '157.37'
$ws.Range("E37").Value = '  -0.01%  '
$ws.Range("D38").Value = '18.87'
$ws.Range("E38").Value = '  +0.25%  '
$ws.Range("E39").Value = '  +1.43%  '
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("D41").Value = '4.90'
$ws.Range("E41").Value = '  +2.63%  '
$ws.Range("D43").Value = '1.63'
$ws.Range("E43").Value = '  -4.15%  '
$ws.Range("D44").Value = '2.38'
$ws.Range("E44").Value = '  -3.88%  '
$ws.Range("E45").Value = '  -11.48%  '
$ws.Range("D46").Value = '38.79'
$ws.Range("E46").Value = '  +0.45%  '
$ws.Range("D47").Value = '146.31'
$ws.Range("E47").Value = '  -0.50%  '
$ws.Range("D48").Value = '0.541'
$ws.Range("E48").Value = '  +1.83%  '
$ws.Range("E49").Value = '  -0.31%  '
$ws.Range("E50").Value = '  -0.68%  '
$ws.Range("D51").Value = '0.0742'
$ws.Range("E51").Value = '  +1.25%  '
